$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Adicionados balanços concatenados em uma única planilha."
# Clear the (all-zero) numeric values on rows 64 and 79, columns C and E:AO,
# leaving them as blank cells (B/D on these rows were already blank).
$ws.Range("C64").ClearContents()
$ws.Range("E64:AO64").ClearContents()

$ws.Range("C79").ClearContents()
$ws.Range("E79:AO79").ClearContents()
